$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text starts with a given prefix, so the
# edits below don't rely on brittle hard-coded paragraph indices.
# ---------------------------------------------------------------------------
function Find-ParagraphByPrefix([string]$prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Helper: wrap a <w:body> fragment + extra relationships into a full
# WordprocessingML "insertable" package and push it into a Range via
# InsertXML. Matching hyperlink targets are de-duplicated against existing
# relationships automatically by the host, so re-used URLs keep their
# original rIds and only genuinely new URLs mint new relationships.
# ---------------------------------------------------------------------------
# Minimal styles part (carried along so rStyle="Hyperlink" / pStyle="NoSpacing"
# survive the InsertXML round-trip instead of being silently dropped because
# they're "unknown" inside the standalone inserted package).
$stylesPart = '<w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/></w:style><w:style w:type="character" w:default="1" w:styleId="DefaultParagraphFont"><w:name w:val="Default Paragraph Font"/></w:style><w:style w:type="character" w:styleId="Hyperlink"><w:name w:val="Hyperlink"/><w:basedOn w:val="DefaultParagraphFont"/><w:uiPriority w:val="99"/><w:unhideWhenUsed/><w:rPr><w:color w:val="0000FF"/><w:u w:val="single"/></w:rPr></w:style><w:style w:type="paragraph" w:styleId="NoSpacing"><w:name w:val="No Spacing"/><w:basedOn w:val="Normal"/><w:uiPriority w:val="1"/><w:qFormat/><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:style></w:styles>'

function New-PkgXml([string]$bodyInner, [string]$relsInner) {
    return @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>$bodyInner<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">$relsInner<Relationship Id="rIdStylesPart" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml">
<pkg:xmlData>$stylesPart</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1) EDFA paragraph: "Programmers manual: <url>" -> hyperlink display text
#    becomes "Programmers manual:" (two runs), trailing space stays plain.
# ---------------------------------------------------------------------------
$body1 = '<w:p><w:hyperlink r:id="rIdH1" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Programmers</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve"> manual:</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$rels1 = '<Relationship Id="rIdH1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.bbnint.co.uk/documents/data_sheets/Lightwaves_2020/EZTEST-.pdf" TargetMode="External"/>'

$p1 = Find-ParagraphByPrefix("Programmers manual: https://www.bbnint.co.uk")
$p1.Range.InsertXML((New-PkgXml $body1 $rels1))

# ---------------------------------------------------------------------------
# 2) Oscilloscope paragraph: re-flow "Programmers manual" and "Tutorial"
#    text into the hyperlink display runs; drop the manual <w:br/>; add a
#    ", " separator; split "requires" across two runs per the diff.
# ---------------------------------------------------------------------------
$body2 = '<w:p><w:r><w:t xml:space="preserve">Model: TBS 1052B-EDU, </w:t></w:r><w:hyperlink r:id="rIdH2" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Programmers</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>m</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>anual</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve">, </w:t></w:r><w:hyperlink r:id="rIdH3" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Tutorial</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">requires </w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:t xml:space="preserve"> party driver </w:t></w:r><w:hyperlink r:id="rIdH4" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>TEKVISA</w:t></w:r></w:hyperlink></w:p>'
$rels2 = '<Relationship Id="rIdH2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.tek.com/oscilloscope/tds1000-manual" TargetMode="External"/><Relationship Id="rIdH3" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.element14.com/community/groups/test-and-measurement/blog/2014/08/22/how-to-connect-a-tektronix-tds2024b-to-a-windows-7-computer-using-usb" TargetMode="External"/><Relationship Id="rIdH4" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.element14.com/community/external-link.jspa?url=http%3A%2F%2Fwww.tek.com%2Foscilloscope%2Ftds7054-software-2" TargetMode="External"/>'

$p2 = Find-ParagraphByPrefix("Model: TBS 1052B-EDU")
$p2.Range.InsertXML((New-PkgXml $body2 $rels2))

# ---------------------------------------------------------------------------
# 3) MATLAB "passing data to python" link: display text becomes
#    "Install engine" (same target URL/relationship).
# ---------------------------------------------------------------------------
$body3 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:hyperlink r:id="rIdH5" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Install engine</w:t></w:r></w:hyperlink></w:p>'
$rels3 = '<Relationship Id="rIdH5" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://au.mathworks.com/help/matlab/matlab_external/passing-data-to-python.html" TargetMode="External"/>'

$p3 = Find-ParagraphByPrefix("https://au.mathworks.com/help/matlab/matlab_external/passing-data-to-python.html")
$p3.Range.InsertXML((New-PkgXml $body3 $rels3))

# ---------------------------------------------------------------------------
# 4) Trailing bookmark-only paragraph becomes three new "NoSpacing"
#    paragraphs, each holding a hyperlink about the MATLAB/python engine.
# ---------------------------------------------------------------------------
$body4 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:hyperlink r:id="rIdH6" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Getting started on engine</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:hyperlink r:id="rIdH7" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Returned data from python</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:hyperlink r:id="rIdH8" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>Pass data to python</w:t></w:r></w:hyperlink></w:p>'
$rels4 = '<Relationship Id="rIdH6" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://au.mathworks.com/help/matlab/matlab_external/get-started-with-matlab-engine-for-python.html" TargetMode="External"/><Relationship Id="rIdH7" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://au.mathworks.com/help/matlab/matlab_external/handle-data-returned-from-python.html" TargetMode="External"/><Relationship Id="rIdH8" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://au.mathworks.com/help/matlab/matlab_external/pass-data-to-python.html" TargetMode="External"/>'

$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p4.Range.InsertXML((New-PkgXml $body4 $rels4))

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
